# Major updates in NA scripts
#
# The contact list is updated from a Michigan/TROY "hostess-test" data set
# to a new Ontario/ETOBICOKE "ahostess-test1" data set, and the blank
# template rows (6-16) are removed from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: hostess ---
$ws.Range("A2").Value = "ahostess-test1"
$ws.Range("C2").Value = "ahostess-test1@test.com"
$ws.Range("D2").Value = "450 KIPLING AVE"
$ws.Range("E2").Value = "ETOBICOKE"
$ws.Range("F2").Value = "Ontario"
$ws.Range("G2").Value = "M8Z 5E1"

# --- Row 3: cohost ---
$ws.Range("A3").Value = "bcohost-test1"
$ws.Range("C3").Value = "bcohost-test1@test.com"
$ws.Range("D3").Value = "450 KIPLING AVE"
$ws.Range("E3").Value = "ETOBICOKE"
$ws.Range("F3").Value = "Ontario"
$ws.Range("G3").Value = "M8Z 5E1"

# --- Row 4: guest1 ---
$ws.Range("A4").Value = "guest1-test1"
$ws.Range("C4").Value = "guest1-test1@test.com"
$ws.Range("D4").Value = "450 KIPLING AVE"
$ws.Range("E4").Value = "ETOBICOKE"
$ws.Range("F4").Value = "Ontario"
$ws.Range("G4").Value = "M8Z 5E1"

# --- Row 5: guest2 ---
$ws.Range("A5").Value = "guest2-test1"
$ws.Range("C5").Value = "guest2-test1@test.com"
$ws.Range("D5").Value = "450 KIPLING AVE"
$ws.Range("E5").Value = "ETOBICOKE"
$ws.Range("F5").Value = "Ontario"
$ws.Range("G5").Value = "M8Z 5E1"

# Remove the now-unused blank template rows 6-16 entirely (shrinks the
# sheet's used range down to A1:G5).
$ws.Range("A6:A16").EntireRow.Delete()

# Match the author's final UI selection state.
$ws.Range("C20").Select()
